# Update the confusion-matrix counts/percentages with the corrected
# translations (the model's predicted vs. true label numbers changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "17 (0.8095)"
$ws.Range("C4").Value = "2 (0.1429)"
$ws.Range("D3").Value = "4 (0.1905)"
$ws.Range("D4").Value = "12 (0.8571)"

# Move the (out-of-grid) selection from G4 to G3, matching the saved view.
$ws.Range("G3").Select() | Out-Null
